# The "bgwrer" company turned out to be a duplicate/bogus entry that was
# creating a second map pin and throwing off the location filters, so we
# drop it from the Companies and Locations lookup tables. We also correct
# the AssetTypes color value, which had been mistyped.

$wb = $excel.ActiveWorkbook

$companies = $wb.Worksheets.Item("Companies")
$companies.Rows.Item(3).Delete()

$locations = $wb.Worksheets.Item("Locations")
$locations.Rows.Item(3).Delete()

$assetTypes = $wb.Worksheets.Item("AssetTypes")
$assetTypes.Range("C2").Value = "#213916"
